# Daily attendance processing - reorders the "Recorded By" audit trail
# (column G) so entries read from most-authoritative (System) first,
# without disturbing rows whose trail is anchored by an admin override
# (admin@admin.com) or rows that only have a single recorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -eq $null) { continue }
    if ($text -eq "") { continue }
    if ($text -notlike "*,*") { continue }
    if ($text -like "*admin@admin.com*") { continue }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = $reversed -join ", "
}
